{"js": "// Apply the \"Last Mash Log\" updates:\n//  - Batch date: 2019-03-29 -> 2019-04-17\n//  - Batch/lot number: 2019/03-4E -> 2019/04-1A\n//  - Corn weight: 5500 -> 500\n//  - Corn order number: 320, 321 -> 16132\n//  - Rye order number: 1232 -> 16110\n//  - Malted Barley weight: 890 -> 1100\n\nconst replacements = [\n  { find: \"2019-03-29\", replace: \"2019-04-17\" },\n  { find: \"2019/03-4E\", replace: \"2019/04-1A\" },\n  { find: \"5500\", replace: \"500\" },\n  { find: \"320, 321\", replace: \"16132\" },\n  { find: \"1232\", replace: \"16110\" },\n  { find: \"890\", replace: \"1100\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  // Replace the first (and expected-only) match.\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the \"Last Mash Log\" updates:\n#  - Batch date: 2019-03-29 -> 2019-04-17\n#  - Batch/lot number: 2019/03-4E -> 2019/04-1A\n#  - Corn weight: 5500 -> 500\n#  - Corn order number: 320, 321 -> 16132\n#  - Rye order number: 1232 -> 16110\n#  - Malted Barley weight: 890 -> 1100\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{ Find = \"2019-03-29\"; Replace = \"2019-04-17\" },\n    @{ Find = \"2019/03-4E\"; Replace = \"2019/04-1A\" },\n    @{ Find = \"5500\";       Replace = \"500\" },\n    @{ Find = \"320, 321\";   Replace = \"16132\" },\n    @{ Find = \"1232\";       Replace = \"16110\" },\n    @{ Find = \"890\";        Replace = \"1100\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($r.Find, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $r.Replace, $wdReplaceAll)\n}\n\n$d.Save()\n"}
